$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.030.94'
$ws.Range("E2").Value = '  -0.59%  '
$ws.Range("D3").Value = '1.753.08'
$ws.Range("E3").Value = '  -3.15%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '336.90'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").Value = '0.9984'
$ws.Range("E6").Value = '  -0.15%  '
$ws.Range("D7").Value = '0.3766'
$ws.Range("E7").Value = '  -4.31%  '
$ws.Range("D8").Value = '0.3351'
$ws.Range("E8").Value = '  -4.38%  '
$ws.Range("D9").Value = '45.35'
$ws.Range("E9").Value = '  -5.71%  '
$ws.Range("E10").Value = '  -5.00%  '
$ws.Range("D11").Value = '0.07192'
$ws.Range("E11").Value = '  -4.75%  '
$ws.Range("D12").Value = '0.9994'
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '22.51'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").Value = '6.145'
$ws.Range("E14").Value = '  -5.66%  '
$ws.Range("D15").Value = '7.133'
$ws.Range("E15").Value = '  -0.46%  '
$ws.Range("D16").Value = '1.757.57'
$ws.Range("E16").Value = '  -2.80%  '
$ws.Range("D17").Value = '0.00001055'
$ws.Range("E17").Value = '  -4.33%  '
$ws.Range("D18").Value = '0.06586'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("D19").Value = '80.43'
$ws.Range("E19").Value = '  -5.67%  '
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").Value = '16.86'
$ws.Range("E21").Value = '  -5.05%  '
$ws.Range("D22").Value = '6.236'
$ws.Range("E22").Value = '  -4.99%  '
$ws.Range("D23").Value = '28.077.48'
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").Value = '11.64'
$ws.Range("E24").Value = '  -6.03%  '
$ws.Range("D25").Value = '2.387'
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("D26").Value = '153.08'
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").Value = '19.82'
$ws.Range("E27").Value = '  -7.92%  '
$ws.Range("D28").Value = '2.316'
$ws.Range("E28").Value = '  -8.39%  '
$ws.Range("D29").Value = '1.955.68'
$ws.Range("E29").Value = '  -2.95%  '
$ws.Range("D30").Value = '131.89'
$ws.Range("E30").Value = '  -3.08%  '
$ws.Range("D31").Value = '1.241'
$ws.Range("E31").Value = '  -16.17%  '
$ws.Range("D32").Value = '4.016'
$ws.Range("D33").Value = '5.777'
$ws.Range("E33").Value = '  -7.24%  '
$ws.Range("D34").Value = '0.08696'
$ws.Range("E34").Value = '  -1.85%  '
$ws.Range("D35").Value = '12.18'
$ws.Range("E35").Value = '  -8.11%  '
$ws.Range("D36").Value = '0.6660'
$ws.Range("E36").Value = '  -3.97%  '
$ws.Range("D37").Value = '0.02320'
$ws.Range("E37").Value = '  -4.84%  '
$ws.Range("D38").Value = '0.06186'
$ws.Range("E38").Value = '  -5.18%  '
$ws.Range("D39").Value = '5.153'
$ws.Range("E39").Value = '  -5.67%  '
$ws.Range("D40").Value = '0.2105'
$ws.Range("E40").Value = '  -4.94%  '
$ws.Range("D41").Value = '1.216'
$ws.Range("E41").Value = '  -3.52%  '
$ws.Range("D42").Value = '1.445'
$ws.Range("E42").Value = '  -10.04%  '
$ws.Range("D43").Value = '8.009'
$ws.Range("E43").Value = '  -6.34%  '
$ws.Range("D44").Value = '0.9993'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '13.67'
$ws.Range("E45").Value = '  -6.65%  '
$ws.Range("D46").Value = '3.831'
$ws.Range("D47").Value = '0.6040'
$ws.Range("E47").Value = '  -6.08%  '
$ws.Range("D48").Value = '128.40'
$ws.Range("E48").Value = '  -2.16%  '
$ws.Range("E49").Value = '  -6.30%  '
$ws.Range("D50").Value = '0.07154'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("D51").Value = '1.170'
$ws.Range("E51").Value = '  +0.79%  '
